$d = $word.ActiveDocument

$d.Content.Find.Execute("822×2=1644", $false, $false, $false, $false, $false, $true, 1, $false, "573×7=4011", 2) | Out-Null
$d.Content.Find.Execute("669×7=4683", $false, $false, $false, $false, $false, $true, 1, $false, "648×5=3240", 2) | Out-Null
$d.Content.Find.Execute("850×2=1700", $false, $false, $false, $false, $false, $true, 1, $false, "301×3=903", 2) | Out-Null
$d.Content.Find.Execute("983×8=7864", $false, $false, $false, $false, $false, $true, 1, $false, "946×9=8514", 2) | Out-Null
$d.Content.Find.Execute("143×4=572", $false, $false, $false, $false, $false, $true, 1, $false, "797×6=4782", 2) | Out-Null
$d.Content.Find.Execute("139×4=556", $false, $false, $false, $false, $false, $true, 1, $false, "255×9=2295", 2) | Out-Null
$d.Content.Find.Execute("213×5=1065", $false, $false, $false, $false, $false, $true, 1, $false, "868×8=6944", 2) | Out-Null
$d.Content.Find.Execute("627×8=5016", $false, $false, $false, $false, $false, $true, 1, $false, "624×3=1872", 2) | Out-Null
$d.Content.Find.Execute("307×2=614", $false, $false, $false, $false, $false, $true, 1, $false, "538×5=2690", 2) | Out-Null
$d.Content.Find.Execute("586×7=4102", $false, $false, $false, $false, $false, $true, 1, $false, "758×3=2274", 2) | Out-Null
$d.Content.Find.Execute("518×6=3108", $false, $false, $false, $false, $false, $true, 1, $false, "518×2=1036", 2) | Out-Null
$d.Content.Find.Execute("371×6=2226", $false, $false, $false, $false, $false, $true, 1, $false, "272×3=816", 2) | Out-Null
$d.Content.Find.Execute("502×3=1506", $false, $false, $false, $false, $false, $true, 1, $false, "706×2=1412", 2) | Out-Null
$d.Content.Find.Execute("268×8=2144", $false, $false, $false, $false, $false, $true, 1, $false, "865×8=6920", 2) | Out-Null
$d.Content.Find.Execute("506×3=1518", $false, $false, $false, $false, $false, $true, 1, $false, "305×2=610", 2) | Out-Null
$d.Content.Find.Execute("245×4=980", $false, $false, $false, $false, $false, $true, 1, $false, "291×7=2037", 2) | Out-Null
$d.Content.Find.Execute("954×2=1908", $false, $false, $false, $false, $false, $true, 1, $false, "917×3=2751", 2) | Out-Null
$d.Content.Find.Execute("988×9=8892", $false, $false, $false, $false, $false, $true, 1, $false, "660×5=3300", 2) | Out-Null
$d.Content.Find.Execute("782×5=3910", $false, $false, $false, $false, $false, $true, 1, $false, "731×2=1462", 2) | Out-Null
$d.Content.Find.Execute("267×3=801", $false, $false, $false, $false, $false, $true, 1, $false, "834×5=4170", 2) | Out-Null
$d.Content.Find.Execute("823×6=4938", $false, $false, $false, $false, $false, $true, 1, $false, "501×4=2004", 2) | Out-Null
$d.Content.Find.Execute("836×7=5852", $false, $false, $false, $false, $false, $true, 1, $false, "194×3=582", 2) | Out-Null
$d.Content.Find.Execute("394×6=2364", $false, $false, $false, $false, $false, $true, 1, $false, "537×9=4833", 2) | Out-Null
$d.Content.Find.Execute("825×9=7425", $false, $false, $false, $false, $false, $true, 1, $false, "287×2=574", 2) | Out-Null
$d.Content.Find.Execute("570×6=3420", $false, $false, $false, $false, $false, $true, 1, $false, "200×8=1600", 2) | Out-Null
